# New crime data collected - update weekly CompStat figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) -----------------------
# "Volume 31   Number  9"  -> "Volume 31   Number  10"
# Only the trailing run ("9") changes; replace just that substring.
$ws.Range("A8").Characters(21, 1).Text = "10"

# "Report Covering the Week  2/26/2024  Through  3/3/2024"
#  -> "Report Covering the Week  3/4/2024  Through  3/10/2024"
# Replace the second (later) date first so the first date's offset is unaffected.
$ws.Range("C9").Characters(47, 8).Text = "3/10/2024"
$ws.Range("C9").Characters(27, 9).Text = "3/4/2024"

# --- Helper "template" cells used to carry number formats across ---------
# (rows 14/15/18/22/27 are not touched by this edit, so they're stable
#  sources for copying the exact existing style indices around)
$textTemplateZero   = $ws.Range("C22")   # style: text, shared string "0"
$textTemplateStar   = $ws.Range("E22")   # style: text, shared string "***.*"
$intTemplate        = $ws.Range("I15")   # style: integer number format
$decTemplate        = $ws.Range("L15")   # style: decimal number format

function Set-AsInt($rangeAddr, $value) {
    $intTemplate.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
    $ws.Range($rangeAddr).Value = $value
}

function Set-AsDecimal($rangeAddr, $value) {
    $decTemplate.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
    $ws.Range($rangeAddr).Value = $value
}

function Set-AsTextZero($rangeAddr) {
    $ws.Range($rangeAddr).NumberFormat = "@"
    $ws.Range($rangeAddr).Value = "0"
    $textTemplateZero.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
}

function Set-AsTextStar($rangeAddr) {
    $ws.Range($rangeAddr).Value = "***.*"
    $textTemplateStar.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4122)
}

# --- Row 16 (Robbery): values only, styles unchanged ----------------------
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = 600
$ws.Range("I16").Value = 10
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 233.333333333333
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = -50

# --- Row 17 (Fel. Assault): several text->number conversions -------------
Set-AsInt "C17" 1
Set-AsInt "D17" 1
Set-AsDecimal "E17" 0
$ws.Range("F17").Value = 2
Set-AsInt "G17" 1
Set-AsDecimal "H17" 100
$ws.Range("I17").Value = 2
Set-AsInt "J17" 1
Set-AsDecimal "K17" 100
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = -75
# M17 unchanged ("***.*")

# --- Row 19 (Gr. Larceny) -------------------------------------------------
Set-AsInt "C19" 1
Set-AsTextZero "D19"
Set-AsTextStar "E19"
$ws.Range("I19").Value = 3
$ws.Range("K19").Value = -25
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -40
$ws.Range("N19").Value = -70

# --- Row 21 (TOTAL): values only, bold styles unchanged -------------------
$ws.Range("C21").Value = 3
$ws.Range("E21").Value = 200
$ws.Range("F21").Value = 10
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 233.333333333333
$ws.Range("I21").Value = 16
$ws.Range("J21").Value = 7
$ws.Range("K21").Value = 128.571428571429
$ws.Range("L21").Value = 60
$ws.Range("M21").Value = 77.777777777777
$ws.Range("N21").Value = -60.975609756097

# --- Row 24 (Petit Larceny) -----------------------------------------------
Set-AsTextZero "C24"
Set-AsInt "D24" 1
Set-AsDecimal "E24" -100
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = -50
$ws.Range("J24").Value = 7
$ws.Range("K24").Value = -85.714285714285

# --- Row 26 (Misd. Assault) -----------------------------------------------
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = -25
$ws.Range("M26").Value = -40

# --- Row 28 (Other Sex Crimes) --------------------------------------------
Set-AsInt "C28" 1
Set-AsInt "F28" 1
Set-AsInt "I28" 1
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
